# Updated cryptos list on Mon Mar  4 15:22:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D whose new value looks like a plain number need to be
# forced to Text format first, otherwise Excel would store them as numeric values
# (losing trailing zeros / producing floating point artifacts) instead of the
# literal text strings used in the source data.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.314.42"
$ws.Range("E2").Value = "  +6.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.558.92"
$ws.Range("E3").Value = "  +3.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.16%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.76"
$ws.Range("E5").Value = "  +1.07%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.11"
$ws.Range("E6").Value = "  +2.06%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.660"
$ws.Range("E7").Value = "  +6.13%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.549.14"
$ws.Range("E8").Value = "  +3.73%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.03%  "

# Row 10 - Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.785"
$ws.Range("E10").Value = "  +8.11%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.169"
$ws.Range("E11").Value = "  +21.62%  "

# Row 12 - now ShibaInu (was Avalanche)
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000288"
$ws.Range("E12").Value = "  +34.35%  "

# Row 13 - now Avalanche (was ShibaInu)
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.40"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.11"
$ws.Range("E14").Value = "  +8.89%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.120.98"
$ws.Range("E15").Value = "  +3.93%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +0.19%  "

# Row 17 - Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.63"
$ws.Range("E17").Value = "  +0.53%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.582.03"
$ws.Range("E18").Value = "  +4.36%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +1.25%  "

# Row 20 - Polygon
$ws.Range("E20").Value = "  +3.15%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "66.267.33"
$ws.Range("E21").Value = "  +6.43%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.37"
$ws.Range("E22").Value = "  -3.06%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.43"
$ws.Range("E23").Value = "  -0.55%  "

# Row 24 - ImmutableX
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.26"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.20"
$ws.Range("E25").Value = "  -2.76%  "

# Row 26 - PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.39"
$ws.Range("E26").Value = "  +2.07%  "

# Row 27 - Filecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -5.33%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.44"
$ws.Range("E28").Value = "  +3.95%  "

# Row 29 - now Cosmos (was Toncoin)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.47"
$ws.Range("E29").Value = "  +4.42%  "

# Row 30 - now Toncoin (was Cosmos)
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").Value = "  +4.96%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +5.14%  "

# Row 32 - RenderToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.33"
$ws.Range("E32").Value = "  -4.68%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  -2.71%  "

# Row 34 - Dai
$ws.Range("E34").Value = "  -0.15%  "

# Row 35 - InjectiveProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.26"
$ws.Range("E35").Value = "  -3.73%  "

# Row 36 - OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.96"
$ws.Range("E36").Value = "  -2.18%  "

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0507"
$ws.Range("E37").Value = "  +3.92%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0740"
$ws.Range("E38").Value = "  +43.40%  "

# Row 39 - Stellar
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.148"
$ws.Range("E39").Value = "  +10.75%  "

# Row 40 - Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +0.63%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42 - WEMIXToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").Value = "  +4.25%  "

# Row 43 - now Monero (was NEARProtocol)
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.09"
$ws.Range("E43").Value = "  +1.86%  "

# Row 44 - now NEARProtocol (was Monero)
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.47"
$ws.Range("E44").Value = "  +3.48%  "

# Row 45 - LidoDAOToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -1.74%  "

# Row 46 - TheGraph
$ws.Range("E46").Value = "  -3.86%  "

# Row 47 - ARBITRUM
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").Value = "  -3.37%  "

# Row 48 - ThetaToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -6.54%  "

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.147"
$ws.Range("E49").Value = "  +6.68%  "

# Row 50 - Celestia
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "15.78"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51 - ApeXProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  +10.58%  "
